# Update quarterly data: shift Trimestre dates from Q2 (01/04) to Q3 (01/07)
# and refresh Valor (D) / Variação % (E) figures for each row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, "01/07/2015", 803, -3.95),
    @(3, "01/07/2015", 2218, -11.32),
    @(4, "01/07/2015", 1700, -14.14),
    @(5, "01/07/2015", 1967, -23.64),
    @(6, "01/07/2015", 2052, 1.74),
    @(7, "01/07/2015", 1674, 15.37),
    @(8, "01/07/2015", 2862, -1.82),
    @(9, "01/07/2015", 3965, -5.19),
    @(10, "01/07/2016", 870, 8.34),
    @(11, "01/07/2016", 2345, 5.73),
    @(12, "01/07/2016", 1734, 2),
    @(13, "01/07/2016", 1790, -9),
    @(14, "01/07/2016", 1838, -10.43),
    @(15, "01/07/2016", 1388, -17.08),
    @(16, "01/07/2016", 2828, -1.19),
    @(17, "01/07/2016", 4469, 12.71),
    @(18, "01/07/2017", 845, -2.87),
    @(19, "01/07/2017", 2038, -13.09),
    @(20, "01/07/2017", 2183, 25.89),
    @(21, "01/07/2017", 1746, -2.46),
    @(22, "01/07/2017", 1668, -9.25),
    @(23, "01/07/2017", 1537, 10.73),
    @(24, "01/07/2017", 2984, 5.52),
    @(25, "01/07/2017", 3943, -11.77),
    @(26, "01/07/2018", 965, 14.2),
    @(27, "01/07/2018", 1878, -7.85),
    @(28, "01/07/2018", 1829, -16.22),
    @(29, "01/07/2018", 1777, 1.78),
    @(30, "01/07/2018", 1948, 16.79),
    @(31, "01/07/2018", 1641, 6.77),
    @(32, "01/07/2018", 2884, -3.35),
    @(33, "01/07/2018", 4007, 1.62),
    @(34, "01/07/2019", $null, -21.35),
    @(35, "01/07/2019", 1883, 0.27),
    @(36, "01/07/2019", 1592, -12.96),
    @(37, "01/07/2019", 1615, -9.119999999999999),
    @(38, "01/07/2019", 2041, 4.77),
    @(39, "01/07/2019", 1314, -19.93),
    @(40, "01/07/2019", 3029, 5.03),
    @(41, "01/07/2019", 3859, -3.69),
    @(42, "01/07/2020", 809, 6.59),
    @(43, "01/07/2020", 1900, 0.9),
    @(44, "01/07/2020", 1625, 2.07),
    @(45, "01/07/2020", 2329, 44.21),
    @(46, "01/07/2020", 1830, -10.34),
    @(47, "01/07/2020", 1762, 34.09),
    @(48, "01/07/2020", 2340, -22.75),
    @(49, "01/07/2020", 3628, -5.99),
    @(50, "01/07/2021", 866, 7.05),
    @(51, "01/07/2021", 2395, 26.05),
    @(52, "01/07/2021", 1375, -15.38),
    @(53, "01/07/2021", 1689, -27.48),
    @(54, "01/07/2021", 1571, -14.15),
    @(55, "01/07/2021", 1275, -27.64),
    @(56, "01/07/2021", 2494, 6.58),
    @(57, "01/07/2021", 4043, 11.44),
    @(58, "01/07/2022", 938, 8.31),
    @(59, "01/07/2022", 1948, -18.66),
    @(60, "01/07/2022", 1429, 3.93),
    @(61, "01/07/2022", 1741, 3.08),
    @(62, "01/07/2022", 1809, 15.15),
    @(63, "01/07/2022", 1233, -3.29),
    @(64, "01/07/2022", 3000, 20.29),
    @(65, "01/07/2022", 3507, -13.26),
    @(66, "01/07/2023", 808, -13.86),
    @(67, "01/07/2023", 2201, 12.99),
    @(68, "01/07/2023", 1503, 5.18),
    @(69, "01/07/2023", 1654, -5),
    @(70, "01/07/2023", 1887, 4.31),
    @(71, "01/07/2023", 1296, 5.11),
    @(72, "01/07/2023", 2751, -8.300000000000001),
    @(73, "01/07/2023", 3686, 5.1),
    @(74, "01/07/2024", 933, 15.47),
    @(75, "01/07/2024", 1831, -16.81),
    @(76, "01/07/2024", 1743, 15.97),
    @(77, "01/07/2024", 2102, 27.09),
    @(78, "01/07/2024", 2067, 9.539999999999999),
    @(79, "01/07/2024", 1334, 2.93),
    @(80, "01/07/2024", 2664, -3.16),
    @(81, "01/07/2024", 4104, 11.34)
)

foreach ($row in $updates) {
    $r = $row[0]
    $newDate = $row[1]
    $newValor = $row[2]
    $newVariacao = $row[3]

    # Column C (Trimestre) is stored as text, not a date - force text format
    # so Excel doesn't auto-convert the "dd/mm/yyyy"-looking string into a
    # date serial number.
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $newDate

    if ($null -ne $newValor) {
        $ws.Cells.Item($r, 4).Value = $newValor
    }

    $ws.Cells.Item($r, 5).Value = $newVariacao
}
